$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd email address that is both the text and the hyperlink
# target of cell A3 ("prunturt@yahoo.fr" -> "pruntrut@yahoo.fr"), while
# keeping the cell's hyperlink style and the visible contents of A1/A2
# untouched.
$ws.Range("A3").Hyperlinks.Delete()
$ws.Range("A3").Value = "pruntrut@yahoo.fr"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:pruntrut@yahoo.fr")
$ws.Range("A3").Style = "Lien hypertexte"

# Restore the selection that was active when the sheet was last saved.
$ws.Range("E7").Select()
